# led-bar-pcb BOM update:
#  - resistor value R1-R10 changes from "10K" text to numeric 51 (51K)
#  - add a note about installing on the backside of the PCB (also referenced by the new connector row)
#  - add a new BOM row for the 2-pin connector (22-28-4112 / Mouser) with its DigiKey-style hyperlink

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: resistor bank value becomes a plain number (was shared text "10K")
$ws.Range("C3").Value = 51

# Row 2: new note cell
$ws.Range("F2").Value = "NOTE - install on backside of PCB"

# Row 4: new Connector line
$ws.Range("A4").Value = "Connector"
$ws.Range("B4").Value = 1

# D4 gets the part number + hyperlink
$ws.Range("D4").Value = "22-28-4112"
$ws.Range("E4").Value = "Mouser"
$ws.Range("F4").Value = "NOTE - install on backside of PCB"

[void]$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.mouser.com/ProductDetail/Molex/22-28-4112", "", "", "22-28-4112")

# Matching the source file, D4 keeps the plain/default cell style rather than
# the blue underlined hyperlink style used by D2/D3, so reset the font after
# the hyperlink is added (Hyperlinks.Add re-applies the hyperlink style).
$ws.Range("D4").Font.Color = 0
$ws.Range("D4").Font.Underline = $false

# Restore the selection the author left the sheet on.
[void]$ws.Range("H10").Select()
